# Applies the "Update COVAX tracker, removing AEFI related metadata" edit described
# in the commit diff: refreshes timestamps, drops 5 AEFI-specific dataElements rows
# (compacting the table from A1:G31 down to A1:G26), and narrows column A.

$wb = $excel.ActiveWorkbook

# --- 1. "Package info" sheet: bump the Created timestamp and derived Identifier ---
$wsInfo = $wb.Worksheets.Item("Package info")
$wsInfo.Range('B6').Value = '2021-02-02T10:31'
$wsInfo.Range('B7').Value = 'COVAX_TRACKER_V1.0.0_DHIS2.34_2021-02-02T10:31'

# --- 2. "dataElements" sheet: compact rows 13-26 with the post-removal content ---
$wsData = $wb.Worksheets.Item("dataElements")

$wsData.Range('A13').Value = 'COVAX - Cardiovascular Disease'
$wsData.Range('B13').Value = 'COVAX Cardiovascular Disease'
$wsData.Range('C13').Value = 'Covax_Comcond_cardi'
$wsData.Range('D13').Value = 'Cardiovascular disease, including hypertension'
$wsData.Range('E13').Value = 'default'
$wsData.Range('F13').Value = '2021-01-31'
$wsData.Range('G13').Value = 'LNHAYF3qdZl'

$wsData.Range('A14').Value = 'COVAX Previously infected with COVID'
$wsData.Range('B14').Value = 'Previously infected with COVID'
$wsData.Range('C14').Value = ''
$wsData.Range('D14').Value = 'Has the patient been infected with COVID-19 within the last 90 days?'
$wsData.Range('E14').Value = 'default'
$wsData.Range('F14').Value = '2021-01-31'
$wsData.Range('G14').Value = 'LOU9t0aR0z7'

$wsData.Range('A15').Value = 'COVAX - Dose Number'
$wsData.Range('B15').Value = 'COVAX - Dose Number'
$wsData.Range('C15').Value = 'dose_number'
$wsData.Range('D15').Value = 'Defines if this is the first vs the second dose, for example'
$wsData.Range('E15').Value = 'default'
$wsData.Range('F15').Value = '2021-01-31'
$wsData.Range('G15').Value = 'LUIsbsm3okG'

$wsData.Range('A16').Value = 'COVAX - AEFIs present'
$wsData.Range('B16').Value = 'COVAX - AEFIs present'
$wsData.Range('C16').Value = 'AEFIs_present'
$wsData.Range('D16').Value = 'Are there any adverse events present after receiving the immunization?'
$wsData.Range('E16').Value = 'default'
$wsData.Range('F16').Value = '2021-02-01'
$wsData.Range('G16').Value = 'm9PgIDAJGlF'

$wsData.Range('A17').Value = 'COVAX - Immunodeficiency'
$wsData.Range('B17').Value = 'COVAX - Immunodeficiency'
$wsData.Range('C17').Value = 'Covax_Comcond_immuno'
$wsData.Range('D17').Value = ''
$wsData.Range('E17').Value = 'default'
$wsData.Range('F17').Value = '2021-01-31'
$wsData.Range('G17').Value = 'MuZ9dMVXpuM'

$wsData.Range('A18').Value = 'COVAX - Multiple products used - Explain'
$wsData.Range('B18').Value = 'Multiple products explanation'
$wsData.Range('C18').Value = ''
$wsData.Range('D18').Value = ''
$wsData.Range('E18').Value = 'default'
$wsData.Range('F18').Value = '2021-01-31'
$wsData.Range('G18').Value = 'OAxinuYFDG6'

$wsData.Range('A19').Value = 'COVAX - Eligibility'
$wsData.Range('B19').Value = 'COVAX - Eligibility'
$wsData.Range('C19').Value = 'eligible'
$wsData.Range('D19').Value = 'Whether the client is eligible for the vaccine'
$wsData.Range('E19').Value = 'default'
$wsData.Range('F19').Value = '2021-01-31'
$wsData.Range('G19').Value = 'OXBnk0mWcOw'

$wsData.Range('A20').Value = 'COVAX Date of remission'
$wsData.Range('B20').Value = 'Date of remission COVID'
$wsData.Range('C20').Value = ''
$wsData.Range('D20').Value = ''
$wsData.Range('E20').Value = 'default'
$wsData.Range('F20').Value = '2021-01-31'
$wsData.Range('G20').Value = 'oZiX5uHkacg'

$wsData.Range('A21').Value = 'COVAX - Diabetes'
$wsData.Range('B21').Value = 'COVAX Diabetes'
$wsData.Range('C21').Value = 'Covax_Comcond_diabetes'
$wsData.Range('D21').Value = ''
$wsData.Range('E21').Value = 'default'
$wsData.Range('F21').Value = '2021-01-31'
$wsData.Range('G21').Value = 'TT1h0vGu5da'

$wsData.Range('A22').Value = 'AEFI_Nonserious adverse event following immunization'
$wsData.Range('B22').Value = 'Non-serious adverse event following immunization'
$wsData.Range('C22').Value = ''
$wsData.Range('D22').Value = ''
$wsData.Range('E22').Value = 'default'
$wsData.Range('F22').Value = '2021-01-31'
$wsData.Range('G22').Value = 'V5CfmqTsKrv'

$wsData.Range('A23').Value = 'COVAX - Neurological/Neuromuscular'
$wsData.Range('B23').Value = 'Neurological/Neuromuscular'
$wsData.Range('C23').Value = 'Comcond_neuro'
$wsData.Range('D23').Value = 'Chronic neurological or neuromuscular disease'
$wsData.Range('E23').Value = 'default'
$wsData.Range('F23').Value = '2021-01-31'
$wsData.Range('G23').Value = 'VCetMtYu1DY'

$wsData.Range('A24').Value = 'COVAX - Malignancy'
$wsData.Range('B24').Value = 'COVAX - Malignancy'
$wsData.Range('C24').Value = 'Covax_Comcond_malig'
$wsData.Range('D24').Value = 'COVAX - Malignancy'
$wsData.Range('E24').Value = 'default'
$wsData.Range('F24').Value = '2021-01-31'
$wsData.Range('G24').Value = 'xVxLMku5DMX'

$wsData.Range('A25').Value = 'COVAX- Batch Number'
$wsData.Range('B25').Value = 'COVAX - Batch Number'
$wsData.Range('C25').Value = ''
$wsData.Range('D25').Value = 'Batch number for this vaccine dose'
$wsData.Range('E25').Value = 'default'
$wsData.Range('F25').Value = '2021-01-31'
$wsData.Range('G25').Value = 'Yp1F4txx8tm'

$wsData.Range('A26').Value = 'COVAX - Dose Expiry Date'
$wsData.Range('B26').Value = 'COVAX - Dose Expiry Date'
$wsData.Range('C26').Value = 'dose_expiry_date'
$wsData.Range('D26').Value = 'Date dose of expiration'
$wsData.Range('E26').Value = 'default'
$wsData.Range('F26').Value = '2021-01-31'
$wsData.Range('G26').Value = 'YTQulAldGOs'

# The 5 AEFI rows that used to occupy 13..31 are now gone; delete the surplus
# trailing rows so the sheet dimension shrinks from A1:G31 to A1:G26
$wsData.Range('A27:G31').EntireRow.Delete()

# Column A no longer needs to fit the long AEFI question text, narrow it
$wsData.Columns.Item(1).ColumnWidth = 54.7109375

# --- 3. "categories" sheet: bump the Last updated date ---
$wsCat = $wb.Worksheets.Item("categories")
$wsCat.Range('B2').Value = '2021-02-02'
